# Updates for 7 May
# Adds a new "5/6/20" data column (BC) to Sheet1, copying the style of the
# preceding "5/5/20" column (BB) and filling in the new day's death counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the whole BB column's formatting into BC first (mirrors the way
# Excel extends formatting when a user fills a new column to the right),
# then overwrite the values with the new day's numbers.
$ws.Range("BB1:BB54").Copy() | Out-Null
$ws.Range("BC1:BC54").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# New date header, mirroring the format of the prior header cell (BB1).
$ws.Range("BC1").Value = " 5/6/20"

# New day's values for each state row (row 2 .. row 54), aligned to the
# same rows as column BB.
$values = @(343,10,426,87,2462,921,2718,193,277,1539,1327,5,17,66,2974,1377,219,164,283,2167,62,1437,4420,4250,485,374,425,16,86,286,111,8572,169,25956,493,31,1225,253,115,3347,99,370,305,29,239,1006,58,52,713,880,51,362,7)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 55).Value = $values[$i]   # column 55 = BC
}

# Select the new day's first data cell, same as a user tabbing one column
# to the right after filling in the new data (matches the saved selection).
$ws.Range("BC2").Select() | Out-Null

$wb.Save()
